# Atualização automática de preços de eletricidade
# Updates row 2 of the Spot_PT sheet with the latest daily/hourly spot prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46059
$ws.Range("B2").Value = 0.01
$ws.Range("C2").Value = -0.02
$ws.Range("D2").Value = -0.12
$ws.Range("E2").Value = -0.15
$ws.Range("F2").Value = -0.14
$ws.Range("G2").Value = -0.06
$ws.Range("H2").Value = 0.03
$ws.Range("I2").Value = 1.69
$ws.Range("J2").Value = 15.22
$ws.Range("K2").Value = 11.37
$ws.Range("L2").Value = 0.39
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0.01
$ws.Range("P2").Value = -0.01
$ws.Range("Q2").Value = -0.01
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 1.89
$ws.Range("T2").Value = 5.23
$ws.Range("U2").Value = 20.79
$ws.Range("V2").Value = 29.94
$ws.Range("W2").Value = 28.76
$ws.Range("X2").Value = 8.16
$ws.Range("Y2").Value = 2.41
$ws.Range("Z2").Value = 5.22
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 17.32
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 29.35
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 13.3
$ws.Range("AG2").Value = "0h-23h"
